$d = $word.ActiveDocument

$replacements = @(
    @("87×21=", "94×46="),
    @("56×72=", "78×30="),
    @("14×22=", "34×41="),
    @("40×75=", "64×89="),
    @("66×15=", "12×88="),
    @("56×35=", "73×55="),
    @("24×67=", "29×61="),
    @("59×13=", "41×87="),
    @("91×34=", "65×60="),
    @("72×51=", "16×93="),
    @("64×54=", "76×17="),
    @("34×65=", "74×92="),
    @("82×32=", "68×99="),
    @("13×79=", "89×51="),
    @("31×86=", "94×35="),
    @("37×14=", "48×11="),
    @("23×93=", "99×62="),
    @("44×75=", "18×23="),
    @("33×47=", "80×98="),
    @("52×50=", "66×48="),
    @("20×78=", "40×81="),
    @("23×87=", "30×89="),
    @("77×55=", "96×55="),
    @("97×22=", "98×21="),
    @("54×67=", "80×14=")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $d.Content.Find.Execute($old, $true, $true, $false, $false, $false, $true, 1, $false, $new, 2)
}
